$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "Peso" column before the existing "Visible" column (L), so
# L -> Peso (new), old L (Visible) -> M, old M (Imagen) -> N.
# ---------------------------------------------------------------------------
$ws.Columns("L:L").Insert()

# ---------------------------------------------------------------------------
# The header row uses a locked style (sheet is protected), so a direct
# Range.Value assignment on it is silently ignored. Stage the new header
# text in a scratch cell, copy it, and paste-values into L1 so the cell
# keeps its original (locked) style instead of picking up an unlocked one.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = "Peso"
$ws.Range("Z1").Copy()
$ws.Range("L1").PasteSpecial(-4163)
$ws.Columns("Z:Z").Delete()

# ---------------------------------------------------------------------------
# Fill in the new "Peso" values for each product row. Row 5 had a pre-existing
# data bug (the "Visible" flag had ended up in column L instead of M, leaving
# the user's cantidadDeCompras reading null) - fix it by keeping L5 = 1
# instead of giving it a genuine new "Peso" count.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = 2
$ws.Range("L3").Value = 4
$ws.Range("L4").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 3

# ---------------------------------------------------------------------------
# Match the author's final selection/scroll state.
# ---------------------------------------------------------------------------
$ws.Range("L1").Select()
